$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1989795918367347
$ws.Range("C2").Value = 0.5535714285714286
$ws.Range("J2").Value = 0.01275510204081633
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.09183673469387756
$ws.Range("B3").Value = 0.008888888888888889
$ws.Range("C3").Value = 0.02222222222222222
$ws.Range("J3").Value = 0.03111111111111111
$ws.Range("P3").Value = 0.76
$ws.Range("S3").Value = 0.1777777777777778
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("P4").Value = 0.6363636363636364
$ws.Range("S4").Value = 0.303030303030303
$ws.Range("B6").Value = 0.05439330543933055
$ws.Range("D6").Value = 0.004184100418410041
$ws.Range("F6").Value = 0.04602510460251046
$ws.Range("J6").Value = 0.301255230125523
$ws.Range("O6").Value = 0.01673640167364017
$ws.Range("Q6").Value = 0.1380753138075314
$ws.Range("R6").Value = 0.07531380753138076
$ws.Range("S6").Value = 0.3640167364016736
$ws.Range("B7").Value = 0.1475409836065574
$ws.Range("D7").Value = 0.004098360655737705
$ws.Range("E7").Value = 0.004098360655737705
$ws.Range("F7").Value = 0.04508196721311476
$ws.Range("J7").Value = 0.1024590163934426
$ws.Range("O7").Value = 0.03278688524590164
$ws.Range("Q7").Value = 0.1229508196721311
$ws.Range("R7").Value = 0.0860655737704918
$ws.Range("S7").Value = 0.4549180327868853
$ws.Range("B8").Value = 0.1216814159292035
$ws.Range("D8").Value = 0.02212389380530973
$ws.Range("F8").Value = 0.06858407079646017
$ws.Range("J8").Value = 0.1128318584070796
$ws.Range("O8").Value = 0.02876106194690265
$ws.Range("Q8").Value = 0.1615044247787611
$ws.Range("R8").Value = 0.07079646017699115
$ws.Range("S8").Value = 0.413716814159292
$ws.Range("B9").Value = 0.1285714285714286
$ws.Range("D9").Value = 0.03571428571428571
$ws.Range("F9").Value = 0.1
$ws.Range("J9").Value = 0.09285714285714286
$ws.Range("O9").Value = 0.02142857142857143
$ws.Range("Q9").Value = 0.15
$ws.Range("R9").Value = 0.05714285714285714
$ws.Range("S9").Value = 0.4142857142857143
$ws.Range("B10").Value = 0.1504424778761062
$ws.Range("D10").Value = 0.01287208366854384
$ws.Range("E10").Value = 0.003218020917135961
$ws.Range("F10").Value = 0.07320997586484312
$ws.Range("J10").Value = 0.1142397425583266
$ws.Range("O10").Value = 0.01769911504424779
$ws.Range("Q10").Value = 0.1617055510860821
$ws.Range("R10").Value = 0.08125502815768303
$ws.Range("S10").Value = 0.3853580048270314
$ws.Range("G11").Value = 0.1662337662337662
$ws.Range("J11").Value = 0.07532467532467532
$ws.Range("K11").Value = 0.2285714285714286
$ws.Range("L11").Value = 0.5064935064935064
$ws.Range("S11").Value = 0.02337662337662338
$ws.Range("G12").Value = 0.7549019607843137
$ws.Range("J12").Value = 0.1764705882352941
$ws.Range("L12").Value = 0.02450980392156863
$ws.Range("S12").Value = 0.04411764705882353
$ws.Range("G13").Value = 0.6590909090909091
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.09090909090909091
$ws.Range("F15").Value = 0.02868852459016394
$ws.Range("H15").Value = 0.1352459016393443
$ws.Range("I15").Value = 0.04508196721311476
$ws.Range("J15").Value = 0.360655737704918
$ws.Range("K15").Value = 0.0860655737704918
$ws.Range("O15").Value = 0.04918032786885246
$ws.Range("S15").Value = 0.2950819672131147
$ws.Range("F16").Value = 0.0242914979757085
$ws.Range("H16").Value = 0.1781376518218623
$ws.Range("I16").Value = 0.05668016194331984
$ws.Range("J16").Value = 0.3684210526315789
$ws.Range("K16").Value = 0.145748987854251
$ws.Range("M16").Value = 0.0242914979757085
$ws.Range("O16").Value = 0.06477732793522267
$ws.Range("S16").Value = 0.1376518218623482
$ws.Range("F17").Value = 0.03380281690140845
$ws.Range("H17").Value = 0.1690140845070423
$ws.Range("I17").Value = 0.08169014084507042
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.1098591549295775
$ws.Range("M17").Value = 0.01126760563380282
$ws.Range("O17").Value = 0.07042253521126761
$ws.Range("S17").Value = 0.123943661971831
$ws.Range("F18").Value = 0.0335195530726257
$ws.Range("H18").Value = 0.1843575418994413
$ws.Range("I18").Value = 0.08379888268156424
$ws.Range("J18").Value = 0.3407821229050279
$ws.Range("K18").Value = 0.106145251396648
$ws.Range("M18").Value = 0.0335195530726257
$ws.Range("O18").Value = 0.0782122905027933
$ws.Range("S18").Value = 0.1396648044692737
$ws.Range("F19").Value = 0.02347762289068232
$ws.Range("H19").Value = 0.2090975788701394
$ws.Range("I19").Value = 0.05429200293470286
$ws.Range("J19").Value = 0.3565663976522377
$ws.Range("K19").Value = 0.1327953044754219
$ws.Range("M19").Value = 0.02347762289068232
$ws.Range("N19").Value = 0.001467351430667645
$ws.Range("O19").Value = 0.07116654438738078
$ws.Range("S19").Value = 0.1276595744680851
